$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")
$ws.Range("A1").Value = "I18n"
